$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.502717018127441
$ws.Range("B1").Value = 3.65783429145813
$ws.Range("C1").Value = 5.986981391906738
$ws.Range("D1").Value = 1.45356810092926
$ws.Range("E1").Value = 0.8501076698303223
